$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1345.4546
$ws.Range("I40").Value = 1080
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 1080
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -905
$ws.Range("N40").Value = -4350
$ws.Range("H52").Value = 15000
$ws.Range("J52").Value = 15000
$ws.Range("L52").Value = 45000
$ws.Range("N52").Value = -45320
$ws.Range("H58").Value = 1229.0625
$ws.Range("I58").Value = 1118.9286
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 3356.7858
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -3206.7858
$ws.Range("N58").Value = -6300
$ws.Range("H64").Value = 3547
$ws.Range("I64").Value = 3317.5
$ws.Range("J64").Value = 3700
$ws.Range("K64").Value = 3317.5
$ws.Range("L64").Value = 3700
$ws.Range("M64").Value = -3069.5
$ws.Range("N64").Value = -4196
$ws.Range("H67").Value = 3547
$ws.Range("I67").Value = 3317.5
$ws.Range("J67").Value = 3700
$ws.Range("K67").Value = 3317.5
$ws.Range("L67").Value = 3700
$ws.Range("M67").Value = -2459.5
$ws.Range("N67").Value = -5416
$ws.Range("H74").Value = 3095.4546
$ws.Range("I74").Value = 3125
$ws.Range("J74").Value = 3088.889
$ws.Range("K74").Value = 3125
$ws.Range("L74").Value = 3088.889
$ws.Range("M74").Value = -2189
$ws.Range("N74").Value = -4960.889
$ws.Range("H76").Value = 3166.6667
$ws.Range("I76").Value = 3166.6667
$ws.Range("K76").Value = 3166.6667
$ws.Range("M76").Value = -2851.6667
$ws.Range("H77").Value = 3095.4546
$ws.Range("I77").Value = 3125
$ws.Range("J77").Value = 3088.889
$ws.Range("K77").Value = 15625
$ws.Range("L77").Value = 15444.445
$ws.Range("M77").Value = -10945
$ws.Range("N77").Value = -24804.445
$ws.Range("H79").Value = 3166.6667
$ws.Range("I79").Value = 3166.6667
$ws.Range("K79").Value = 3166.6667
$ws.Range("M79").Value = -2074.6667
$ws.Range("H107").Value = 50000628
$ws.Range("I107").Value = 659.4211
$ws.Range("K107").Value = 659.4211
$ws.Range("M107").Value = 1260.5789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2123.353
$ws.Range("I88").Value = 2958
$ws.Range("J88").Value = 1775.5834
$ws.Range("K88").Value = 2958
$ws.Range("L88").Value = 1775.5834
$ws.Range("M88").Value = -2552
$ws.Range("N88").Value = -2587.5834
$ws.Range("H91").Value = 2123.353
$ws.Range("I91").Value = 2958
$ws.Range("J91").Value = 1775.5834
$ws.Range("K91").Value = 2958
$ws.Range("L91").Value = 1775.5834
$ws.Range("M91").Value = -1554
$ws.Range("N91").Value = -4583.5834
$ws.Range("H132").Value = 2251.8965
$ws.Range("I132").Value = 1959.5333
$ws.Range("J132").Value = 2565.1428
$ws.Range("K132").Value = 5878.5999
$ws.Range("L132").Value = 7695.428400000001
$ws.Range("M132").Value = -3348.5999
$ws.Range("N132").Value = -12755.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1761.6666
$ws.Range("I86").Value = 1610
$ws.Range("J86").Value = 2216.6667
$ws.Range("K86").Value = 1610
$ws.Range("L86").Value = 2216.6667
$ws.Range("M86").Value = -487
$ws.Range("N86").Value = -4462.6667
$ws.Range("H89").Value = 1761.6666
$ws.Range("I89").Value = 1610
$ws.Range("J89").Value = 2216.6667
$ws.Range("K89").Value = 8050
$ws.Range("L89").Value = 11083.3335
$ws.Range("M89").Value = -2434
$ws.Range("N89").Value = -22315.3335
$ws.Range("H105").Value = 2343.5078
$ws.Range("I105").Value = 2350.4644
$ws.Range("K105").Value = 2350.4644
$ws.Range("M105").Value = -603.4643999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2462.5
$ws.Range("I62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("H65").Value = 2462.5
$ws.Range("I65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("H94").Value = 2595.9092
$ws.Range("I94").Value = 4541.3335
$ws.Range("J94").Value = 1866.375
$ws.Range("K94").Value = 4541.3335
$ws.Range("L94").Value = 1866.375
$ws.Range("M94").Value = -4090.3335
$ws.Range("N94").Value = -2768.375
$ws.Range("H134").Value = 3310.8262
$ws.Range("I134").Value = 2792.1052
$ws.Range("J134").Value = 5774.75
$ws.Range("K134").Value = 8376.3156
$ws.Range("L134").Value = 17324.25
$ws.Range("M134").Value = -5841.3156
$ws.Range("N134").Value = -22394.25
$ws.Range("H141").Value = 42741.09
$ws.Range("J141").Value = 42741.09
$ws.Range("L141").Value = 42741.09
$ws.Range("N141").Value = -53101.09

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 521
$ws.Range("I92").Value = 595.1429000000001
$ws.Range("J92").Value = 477.75
$ws.Range("K92").Value = 1785.4287
$ws.Range("L92").Value = 1433.25
$ws.Range("M92").Value = -537.4287000000002
$ws.Range("N92").Value = -3929.25
$ws.Range("H132").Value = 927810.2
$ws.Range("I132").Value = 1948.4736
$ws.Range("J132").Value = 3126731.8
$ws.Range("K132").Value = 17536.2624
$ws.Range("L132").Value = 28140586.2
$ws.Range("M132").Value = -15006.2624
$ws.Range("N132").Value = -28145646.2
$ws.Range("H140").Value = 907.1739
$ws.Range("I140").Value = 683.2222
$ws.Range("J140").Value = 1713.4
$ws.Range("K140").Value = 2049.6666
$ws.Range("L140").Value = 5140.200000000001
$ws.Range("M140").Value = 3130.3334
$ws.Range("N140").Value = -15500.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5394.706
$ws.Range("I70").Value = 5133.3335
$ws.Range("K70").Value = 5133.3335
$ws.Range("M70").Value = -4863.3335
$ws.Range("H73").Value = 5394.706
$ws.Range("I73").Value = 5133.3335
$ws.Range("K73").Value = 5133.3335
$ws.Range("M73").Value = -4197.3335
$ws.Range("H80").Value = 2541.6667
$ws.Range("I80").Value = 2412.5
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 2412.5
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -1414.5
$ws.Range("N80").Value = -4796
$ws.Range("H83").Value = 2541.6667
$ws.Range("I83").Value = 2412.5
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 12062.5
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -7070.5
$ws.Range("N83").Value = -23984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 37780.15
$ws.Range("I22").Value = 111538.89
$ws.Range("J22").Value = 900.7778
$ws.Range("K22").Value = 111538.89
$ws.Range("L22").Value = 900.7778
$ws.Range("M22").Value = -111243.89
$ws.Range("N22").Value = -1490.7778
$ws.Range("H27").Value = 37780.15
$ws.Range("I27").Value = 111538.89
$ws.Range("J27").Value = 900.7778
$ws.Range("K27").Value = 111538.89
$ws.Range("L27").Value = 900.7778
$ws.Range("M27").Value = -111431.89
$ws.Range("N27").Value = -1114.7778
